$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 719, pushing the existing rows 719-747 down to 721-749.
$ws.Rows("719:720").Insert()

# --- New row 719: Acelga, Vega Central Mapocho de Santiago, Primera ---
$ws.Cells.Item(719, 1).Value = 9
$ws.Cells.Item(719, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(719, 3).Value = "Metropolitana"
$ws.Cells.Item(719, 4).Value = 44939
$ws.Cells.Item(719, 5).Value = 13
$ws.Cells.Item(719, 6).Value = 100112009
$ws.Cells.Item(719, 7).Value = "Acelga"
$ws.Cells.Item(719, 8).Value = "Sin especificar"
$ws.Cells.Item(719, 9).Value = "Primera"
$ws.Cells.Item(719, 10).Value = 70
$ws.Cells.Item(719, 11).Value = 12000
$ws.Cells.Item(719, 12).Value = 12000
$ws.Cells.Item(719, 13).Value = 12000
$ws.Cells.Item(719, 14).Value = "`$/docena de atados"
$ws.Cells.Item(719, 15).Value = "Región Metropolitana"
$ws.Cells.Item(719, 16).Value = 4000
$ws.Cells.Item(719, 17).Value = 3
$ws.Cells.Item(719, 18).Value = "Hortaliza"

# --- New row 720: Acelga, Vega Central Mapocho de Santiago, Segunda ---
$ws.Cells.Item(720, 1).Value = 9
$ws.Cells.Item(720, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(720, 3).Value = "Metropolitana"
$ws.Cells.Item(720, 4).Value = 44939
$ws.Cells.Item(720, 5).Value = 13
$ws.Cells.Item(720, 6).Value = 100112009
$ws.Cells.Item(720, 7).Value = "Acelga"
$ws.Cells.Item(720, 8).Value = "Sin especificar"
$ws.Cells.Item(720, 9).Value = "Segunda"
$ws.Cells.Item(720, 10).Value = 52
$ws.Cells.Item(720, 11).Value = 9000
$ws.Cells.Item(720, 12).Value = 9000
$ws.Cells.Item(720, 13).Value = 9000
$ws.Cells.Item(720, 14).Value = "`$/docena de atados"
$ws.Cells.Item(720, 15).Value = "Región Metropolitana"
$ws.Cells.Item(720, 16).Value = 3000
$ws.Cells.Item(720, 17).Value = 3
$ws.Cells.Item(720, 18).Value = "Hortaliza"
